# daschland_ontology/daschland (Alice in DaSCHland)/properties.xlsx
# Commit: "change isPartOfBook to isPartOfBookChapter (#63)"
#
# Row 26 holds the "isPartOfBook" property definition. Rename it (and its
# translated labels / comments / linked-resource class) to describe a link
# to a book *chapter* instead of a whole book.
#
# Columns on Sheet1:
#   A name | B label_en | C label_de | D label_fr | E label_it | F label_rm
#   G comment_en | H comment_de | I comment_fr | J comment_it | K comment_rm
#   L super | M object | N gui_element | O gui_attributes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

# name
$ws.Cells.Item($row, 1).Value2  = "isPartOfBookChapter"

# labels (order mirrors how the author actually typed them: name, label_fr,
# label_de, label_en, label_it)
$ws.Cells.Item($row, 4).Value2  = "Appartient au chapitre"      # label_fr
$ws.Cells.Item($row, 3).Value2  = "Teil des Kapitels"           # label_de
$ws.Cells.Item($row, 2).Value2  = "Part of chapter"             # label_en
$ws.Cells.Item($row, 5).Value2  = "Appartiene al capitolo"      # label_it

# comments
$ws.Cells.Item($row, 7).Value2  = "Belongs to following chapter"          # comment_en
$ws.Cells.Item($row, 8).Value2  = "Gehört zu folgendem Kapitel"           # comment_de
$ws.Cells.Item($row, 9).Value2  = "Appartient au chapitre suivant"        # comment_fr
$ws.Cells.Item($row, 10).Value2 = "Appartiene al seguente capitolo"       # comment_it

# object: the linked resource class changes from :Book to :BookChapter
$ws.Cells.Item($row, 13).Value2 = ":BookChapter"

# leave the author's selection on the row they edited
$ws.Range("A26").Select() | Out-Null
